$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Treino confusion matrix row 9: A9 0 -> 1, B9 19094 -> 19093
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 19093

# Acuracia (Treino)
$ws.Range("B35").Value = 0.9999430459050006

# Precision (Treino)
$ws.Range("B39").Value = 0.9999476275269719

# Recall (Treino)
$ws.Range("B43").Value = 0.9999476275269719

# F1-score (Treino)
$ws.Range("B47").Value = 0.9999476275269719

# Tempo de execucao
$ws.Range("B49").Value = "0:05:21.054732"
